# Update transaction sheet: lowercase headers, fix a couple of dates,
# and drop the scratch WEEKDAY helper formula in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row -> lowercase labels (column order unchanged: date, asset, quantity, price, narration)
$ws.Range("B1").Value = "asset"
$ws.Range("A1").Value = "date"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "price"
$ws.Range("E1").Value = "narration"

# Correct a couple of transaction dates
$ws.Range("A3").Value = 42744
$ws.Range("A5").Value = 42748

# Remove the scratch WEEKDAY(A2) helper formula that lived in I2
$ws.Range("I2").ClearContents()

# Move the selection back to A1 so no stray selection/extra dimension lingers
$ws.Range("A1").Select()
